# Auto-generated edit script: updates market-price derived columns (H-N)
# on multiple rows across several worksheets, per the target diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (index 1) ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = $null
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = $null
$ws.Range("H17").Value = 424.48215
$ws.Range("J17").Value = 430.94446
$ws.Range("L17").Value = 1292.83338
$ws.Range("N17").Value = -1628.83338
$ws.Range("H18").Value = 1102.7142
$ws.Range("I18").Value = 1119.8334
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1119.8334
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -835.8334
$ws.Range("N18").Value = -1568
$ws.Range("H39").Value = 41.63158
$ws.Range("J39").Value = 199
$ws.Range("L39").Value = 597
$ws.Range("N39").Value = -1189
$ws.Range("H47").Value = 77992.25
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 77992.25
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 77992.25
$ws.Range("M47").Value = $null
$ws.Range("N47").Value = -79936.25
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("H86").Value = 9092975
$ws.Range("J86").Value = 2486.8
$ws.Range("L86").Value = 2486.8
$ws.Range("N86").Value = -4732.8
$ws.Range("H89").Value = 9092975
$ws.Range("J89").Value = 2486.8
$ws.Range("L89").Value = 12434
$ws.Range("N89").Value = -23666
$ws.Range("H134").Value = 87876.11
$ws.Range("J134").Value = 87876.11
$ws.Range("L134").Value = 87876.11
$ws.Range("N134").Value = -98016.11
$ws.Range("H136").Value = 78105
$ws.Range("J136").Value = 78105
$ws.Range("L136").Value = 78105
$ws.Range("N136").Value = -88305
$ws.Range("H138").Value = 2828.652
$ws.Range("J138").Value = 3336.6428
$ws.Range("L138").Value = 10009.9284
$ws.Range("N138").Value = -20289.9284

# ---- Sheet: ARM (index 2) ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H6").Value = 238181.81
$ws.Range("I6").Value = 257000
$ws.Range("K6").Value = 257000
$ws.Range("M6").Value = -256827
$ws.Range("H108").Value = 88996.664
$ws.Range("J108").Value = 88996.664
$ws.Range("L108").Value = 88996.664
$ws.Range("N108").Value = -96676.664
$ws.Range("H117").Value = 74584.86
$ws.Range("J117").Value = 74584.86
$ws.Range("L117").Value = 74584.86
$ws.Range("N117").Value = -83762.86
$ws.Range("H118").Value = 49311.332
$ws.Range("J118").Value = 49311.332
$ws.Range("L118").Value = 49311.332
$ws.Range("N118").Value = -52625.332
$ws.Range("H122").Value = 2831
$ws.Range("I122").Value = 2831
$ws.Range("K122").Value = 8493
$ws.Range("M122").Value = -6043
$ws.Range("H137").Value = 110000
$ws.Range("J137").Value = 110000
$ws.Range("L137").Value = 110000
$ws.Range("N137").Value = -120200

# ---- Sheet: BSM (index 3) ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H13").Value = 38085.25
$ws.Range("J13").Value = 38085.25
$ws.Range("L13").Value = 38085.25
$ws.Range("N13").Value = -38421.25
$ws.Range("H26").Value = 18880.666
$ws.Range("I26").Value = 18880.666
$ws.Range("K26").Value = 18880.666
$ws.Range("M26").Value = -18588.666
$ws.Range("H50").Value = 53160.2
$ws.Range("J50").Value = 53160.2
$ws.Range("L50").Value = 53160.2
$ws.Range("N50").Value = -54308.2
$ws.Range("H53").Value = 28664.25
$ws.Range("J53").Value = 28664.25
$ws.Range("L53").Value = 28664.25
$ws.Range("N53").Value = -29812.25
$ws.Range("H109").Value = 77854.14
$ws.Range("J109").Value = 77854.14
$ws.Range("L109").Value = 77854.14
$ws.Range("N109").Value = -80628.14
$ws.Range("H115").Value = 79997.14
$ws.Range("J115").Value = 79997.14
$ws.Range("L115").Value = 79997.14
$ws.Range("N115").Value = -83131.14
$ws.Range("H119").Value = 99995
$ws.Range("J119").Value = 99995
$ws.Range("L119").Value = 99995
$ws.Range("N119").Value = -109671
$ws.Range("H122").Value = 81660
$ws.Range("J122").Value = 81660
$ws.Range("L122").Value = 81660
$ws.Range("N122").Value = -91460
$ws.Range("H135").Value = 82854.28999999999
$ws.Range("J135").Value = 82854.28999999999
$ws.Range("L135").Value = 82854.28999999999
$ws.Range("N135").Value = -92994.28999999999
$ws.Range("H138").Value = 69997.14
$ws.Range("J138").Value = 69997.14
$ws.Range("L138").Value = 69997.14
$ws.Range("N138").Value = -80277.14

# ---- Sheet: CRP (index 4) ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H18").Value = 23328.666
$ws.Range("J18").Value = 23328.666
$ws.Range("L18").Value = 23328.666
$ws.Range("N18").Value = -23788.666
$ws.Range("H31").Value = 2179.2974
$ws.Range("I31").Value = 1555.6364
$ws.Range("K31").Value = 1555.6364
$ws.Range("M31").Value = -1260.6364
$ws.Range("H34").Value = 2179.2974
$ws.Range("I34").Value = 1555.6364
$ws.Range("K34").Value = 1555.6364
$ws.Range("M34").Value = -1353.6364
$ws.Range("H100").Value = 61666.332
$ws.Range("J100").Value = 61666.332
$ws.Range("L100").Value = 61666.332
$ws.Range("N100").Value = -63830.332
$ws.Range("H114").Value = 66201.42999999999
$ws.Range("J114").Value = 73068.336
$ws.Range("L114").Value = 73068.336
$ws.Range("N114").Value = -81746.336
$ws.Range("H133").Value = 74960
$ws.Range("J133").Value = 69940
$ws.Range("L133").Value = 69940
$ws.Range("N133").Value = -75000
$ws.Range("H134").Value = 4204602.5
$ws.Range("I134").Value = 4204602.5
$ws.Range("K134").Value = 12613807.5
$ws.Range("M134").Value = -12611272.5
$ws.Range("H138").Value = 99995
$ws.Range("J138").Value = 99995
$ws.Range("L138").Value = 99995
$ws.Range("N138").Value = -110275

# ---- Sheet: GSM (index 6) ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H46").Value = 26874.75
$ws.Range("J46").Value = 49999.5
$ws.Range("L46").Value = 49999.5
$ws.Range("N46").Value = -50311.5
$ws.Range("H88").Value = 65000
$ws.Range("I88").Value = 65000
$ws.Range("K88").Value = 65000
$ws.Range("M88").Value = -64549
$ws.Range("H91").Value = 65000
$ws.Range("I91").Value = 65000
$ws.Range("K91").Value = 65000
$ws.Range("M91").Value = -63440
$ws.Range("H109").Value = 36016.312
$ws.Range("J109").Value = 36016.312
$ws.Range("L109").Value = 36016.312
$ws.Range("N109").Value = -38096.312
$ws.Range("H110").Value = 99964.836
$ws.Range("J110").Value = 99964.836
$ws.Range("L110").Value = 99964.836
$ws.Range("N110").Value = -108144.836
$ws.Range("H116").Value = 99990
$ws.Range("J116").Value = 99990
$ws.Range("L116").Value = 99990
$ws.Range("N116").Value = -109168
$ws.Range("H119").Value = 87996.664
$ws.Range("J119").Value = 87996.664
$ws.Range("L119").Value = 87996.664
$ws.Range("N119").Value = -97672.664
$ws.Range("H135").Value = 99990.07000000001
$ws.Range("J135").Value = 99990.07000000001
$ws.Range("L135").Value = 99990.07000000001
$ws.Range("N135").Value = -110130.07
$ws.Range("H140").Value = 92138.27
$ws.Range("J140").Value = 92138.27
$ws.Range("L140").Value = 92138.27
$ws.Range("N140").Value = -102498.27

# ---- Sheet: LTW (index 7) ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H42").Value = 15000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 15000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = $null
$ws.Range("N42").Value = -16126
$ws.Range("H49").Value = 15000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 15000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = $null
$ws.Range("N49").Value = -15294
$ws.Range("H117").Value = 41284.777
$ws.Range("J117").Value = 41284.777
$ws.Range("L117").Value = 41284.777
$ws.Range("N117").Value = -50462.777
$ws.Range("H123").Value = 52815.5
$ws.Range("J123").Value = 53300.6
$ws.Range("L123").Value = 53300.6
$ws.Range("N123").Value = -63100.6

# ---- Sheet: WVR (index 8) ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 5225
$ws.Range("I81").Value = 5333.3335
$ws.Range("J81").Value = 4900
$ws.Range("K81").Value = 10666.667
$ws.Range("L81").Value = 9800
$ws.Range("M81").Value = -9605.666999999999
$ws.Range("N81").Value = -11922
$ws.Range("H84").Value = 5225
$ws.Range("I84").Value = 5333.3335
$ws.Range("J84").Value = 4900
$ws.Range("K84").Value = 53333.335
$ws.Range("L84").Value = 49000
$ws.Range("M84").Value = -48029.335
$ws.Range("N84").Value = -59608
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null
